$d = $word.ActiveDocument
